$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.71939999999999
$ws.Range("D12").Value = -8.3705
$ws.Range("E12").Value = 12.45909999999999
$ws.Range("E14").Value = 13.91500000000001
$ws.Range("E22").Value = 11.5255
